$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Logements")
$ws1.Cells.Validation.Delete()
